$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preventing Excel from
# auto-coercing numeric-looking strings (e.g. "566.80") into numbers.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue $ws.Range("D2") '60.857.05'
$ws.Range("E2").Value = '  -3.60%  '
Set-TextValue $ws.Range("D3") '3.358.40'
$ws.Range("E3").Value = '  -2.93%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") '566.80'
$ws.Range("E5").Value = '  -2.14%  '
Set-TextValue $ws.Range("D6") '148.45'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  -1.79%  '
Set-TextValue $ws.Range("D11") '0.412'
$ws.Range("E11").Value = '  +0.94%  '
Set-TextValue $ws.Range("D12") '3.933.91'
$ws.Range("E12").Value = '  -2.91%  '
Set-TextValue $ws.Range("D14") '27.99'
$ws.Range("E14").Value = '  -2.12%  '
Set-TextValue $ws.Range("D15") '3.366.06'
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("E16").Value = '  -1.93%  '
Set-TextValue $ws.Range("D17") '60.968.83'
$ws.Range("E17").Value = '  -3.47%  '
Set-TextValue $ws.Range("D18") '6.36'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("E20").Value = '  -3.49%  '
Set-TextValue $ws.Range("D21") '374.29'
$ws.Range("E21").Value = '  -3.87%  '
Set-TextValue $ws.Range("D22") '75.41'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("E24").Value = '  -0.01%  '
Set-TextValue $ws.Range("D25") '3.503.51'
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("E26").Value = '  -5.52%  '
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("E28").Value = '  -0.08%  '
Set-TextValue $ws.Range("D29") '7.36'
$ws.Range("E29").Value = '  -4.03%  '
$ws.Range("E30").Value = '  +0.00%  '
Set-TextValue $ws.Range("D31") '2.07'
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("E32").Value = '  -5.08%  '
Set-TextValue $ws.Range("D33") '22.86'
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("E34").Value = '  -3.65%  '
$ws.Range("E35").Value = '  +0.44%  '
Set-TextValue $ws.Range("D36") '168.56'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("E37").Value = '  -5.49%  '
$ws.Range("E38").Value = '  -3.98%  '
Set-TextValue $ws.Range("D39") '29.23'
$ws.Range("E39").Value = '  -8.86%  '
Set-TextValue $ws.Range("D40") '3.395.89'
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("E41").Value = '  -3.29%  '
Set-TextValue $ws.Range("D42") '42.29'
$ws.Range("E42").Value = '  -1.28%  '
Set-TextValue $ws.Range("D43") '0.760'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("E44").Value = '  -1.76%  '
$ws.Range("E45").Value = '  -4.26%  '
$ws.Range("E46").Value = '  -6.19%  '
Set-TextValue $ws.Range("D47") '2.497.51'
$ws.Range("E47").Value = '  -3.49%  '
Set-TextValue $ws.Range("D48") '22.62'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("E51").Value = '  -2.60%  '
